$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("overall_results_list")
$ws.Activate()

# Extend the existing "Table1" structured table with a new 4th column.
$lo = $ws.ListObjects.Item("Table1")
$newCol = $lo.ListColumns.Add()

# The best-performing run (row 6, "absolute_error" / "best") gets a note first,
# then the new column gets its header.
$ws.Range("D6").Value = "got best result"
$ws.Range("D1").Value = "Result Desc"

# Give the new column a sensible width.
$ws.Columns.Item(4).ColumnWidth = 14.63

# Widen the visible sheet-tab strip slightly (cosmetic, matches the saved view state).
$excel.ActiveWindow.TabRatio = 0.806

$ws.Range("D2").Select()
